$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 531.8461
$ws.Range("I33").Value = 160.36363
$ws.Range("K33").Value = 160.36363
$ws.Range("M33").Value = 68.63637

# Row 46
$ws.Range("H46").Value = 7209.6772
$ws.Range("I46").Value = 6946.4287
$ws.Range("J46").Value = 9666.666999999999
$ws.Range("K46").Value = 20839.2861
$ws.Range("L46").Value = 29000.001
$ws.Range("M46").Value = -20720.2861
$ws.Range("N46").Value = -29238.001

# Row 60
$ws.Range("H60").Value = 7209.6772
$ws.Range("I60").Value = 6946.4287
$ws.Range("J60").Value = 9666.666999999999
$ws.Range("K60").Value = 20839.2861
$ws.Range("L60").Value = 29000.001
$ws.Range("M60").Value = -20355.2861
$ws.Range("N60").Value = -29968.001

# Row 76
$ws.Range("H76").Value = 3333
$ws.Range("I76").Value = 3333
$ws.Range("K76").Value = 3333
$ws.Range("M76").Value = -3018

# Row 79
$ws.Range("H79").Value = 3333
$ws.Range("I79").Value = 3333
$ws.Range("K79").Value = 3333
$ws.Range("M79").Value = -2241

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 142.85715
$ws.Range("I5").Value = 142.85715
$ws.Range("K5").Value = 142.85715
$ws.Range("M5").Value = -30.85714999999999

# Row 132
$ws.Range("H132").Value = 1524
$ws.Range("I132").Value = 1524
$ws.Range("K132").Value = 4572
$ws.Range("M132").Value = -2042

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 142.85715
$ws.Range("I4").Value = 142.85715
$ws.Range("K4").Value = 142.85715
$ws.Range("M4").Value = -27.85714999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 27785576
$ws.Range("I62").Value = 35721600
$ws.Range("J62").Value = 9499.5
$ws.Range("K62").Value = 35721600
$ws.Range("L62").Value = 9499.5
$ws.Range("M62").Value = -35720976
$ws.Range("N62").Value = -10747.5

# Row 65
$ws.Range("H65").Value = 27785576
$ws.Range("I65").Value = 35721600
$ws.Range("J65").Value = 9499.5
$ws.Range("K65").Value = 178608000
$ws.Range("L65").Value = 47497.5
$ws.Range("M65").Value = -178604880
$ws.Range("N65").Value = -53737.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 353345.25
$ws.Range("I4").Value = 385055.06
$ws.Range("K4").Value = 1155165.18
$ws.Range("M4").Value = -1155053.18

# Row 64
$ws.Range("H64").Value = 994.5
$ws.Range("I64").Value = 994.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2983.5
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("M64").Value = -2713.5

# Row 67
$ws.Range("H67").Value = 994.5
$ws.Range("I67").Value = 994.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2983.5
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("M67").Value = -2047.5

# Row 69
$ws.Range("H69").Value = 13182.4
$ws.Range("I69").Value = 10911.5
$ws.Range("J69").Value = 13750.125
$ws.Range("K69").Value = 32734.5
$ws.Range("L69").Value = 41250.375
$ws.Range("M69").Value = -31923.5
$ws.Range("N69").Value = -42872.375

# Row 72
$ws.Range("H72").Value = 13182.4
$ws.Range("I72").Value = 10911.5
$ws.Range("J72").Value = 13750.125
$ws.Range("K72").Value = 98203.5
$ws.Range("L72").Value = 123751.125
$ws.Range("M72").Value = -94147.5
$ws.Range("N72").Value = -131863.125

# Row 99
$ws.Range("H99").Value = 1375
$ws.Range("I99").Value = 1375
$ws.Range("K99").Value = 4125
$ws.Range("M99").Value = -1879

# Row 127
$ws.Range("H127").Value = 501516.5
$ws.Range("J127").Value = 501516.5
$ws.Range("L127").Value = 1504549.5
$ws.Range("N127").Value = -1514469.5

# Row 139
$ws.Range("H139").Value = 2582.25
$ws.Range("I139").Value = 2276.3333
$ws.Range("K139").Value = 6828.999899999999
$ws.Range("M139").Value = -1688.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 9573.143
$ws.Range("I122").Value = 9573.143
$ws.Range("K122").Value = 28719.429
$ws.Range("M122").Value = -26269.429

# Row 132
$ws.Range("H132").Value = 4674.154
$ws.Range("I132").Value = 4063.4443
$ws.Range("K132").Value = 12190.3329
$ws.Range("M132").Value = -9660.332900000001

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 10751.25
$ws.Range("I20").Value = 10005
$ws.Range("J20").Value = 11000
$ws.Range("K20").Value = 10005
$ws.Range("L20").Value = 11000
$ws.Range("M20").Value = -9779
$ws.Range("N20").Value = -11452

# Row 38
$ws.Range("H38").Value = 23343.334
$ws.Range("I38").Value = 20015
$ws.Range("K38").Value = 20015
$ws.Range("M38").Value = -19605

# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = 0

# Row 100
$ws.Range("H100").Value = 5537.75
$ws.Range("I100").Value = 3667.3333
$ws.Range("J100").Value = 6660
$ws.Range("K100").Value = 3667.3333
$ws.Range("L100").Value = 6660
$ws.Range("M100").Value = -3126.3333
$ws.Range("N100").Value = -7742

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 6667333.5
$ws.Range("I5").Value = 2001
$ws.Range("K5").Value = 2001
$ws.Range("M5").Value = -1889

# Row 80
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -61996

# Row 82
$ws.Range("H82").Value = 45000
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 45000
$ws.Range("N82").Value = -45766

# Row 83
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -189984

# Row 85
$ws.Range("H85").Value = 45000
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 45000
$ws.Range("N85").Value = -47652

# Row 110
$ws.Range("H110").Value = 45622
$ws.Range("J110").Value = 45622
$ws.Range("L110").Value = 45622
$ws.Range("N110").Value = -53802

# Row 126
$ws.Range("H126").Value = 3617.1667
$ws.Range("I126").Value = 3617.1667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10851.5001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -8381.500100000001
